# Adds two new passive skills (Capital City Building Requests / Capital City
# Unit Requests) to the Kingdom Passive Skills sheet, and introduces two new
# "travel time reduction" columns (H & I) used by these new skills.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns at H:I -----------------------------------------
# This shifts the existing effect_type/parent_skill_id/unlocks_at_level/
# is_locked/is_parent columns (old H:L) to J:N, keeping their data intact.
$ws.Range("H1:I1").EntireColumn.Insert()

# --- New column headers -----------------------------------------------------
$ws.Range("H1").Value = "capital_city_building_request_travel_time_reduction"
$ws.Range("I1").Value = "capital_city_unit_request_travel_time_reduction"

# --- New row 19: Capital City Building Requests -----------------------------
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Capital City Building Requests"
$ws.Range("C19").Value = "The requests for repairing or upgrading buildings send out as requests from a capital city will move 15% faster per level for a maximum of a 60% reduction."
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 3
$ws.Range("H19").Value = 0.15
$ws.Range("J19").Value = 12
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 2
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 1

# --- New row 20: Capital City Unit Requests ---------------------------------
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Capital City Unit Requests"
$ws.Range("C20").Value = "When using a capital city to request units be recruited, the time required to travel will be reduced by 15% per level for a max of 60% time reduction."
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 3
$ws.Range("I20").Value = 0.15
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 4
$ws.Range("L20").Value = 2
$ws.Range("M20").Value = 1
$ws.Range("N20").Value = 1

# --- Column widths -----------------------------------------------------------
$ws.Range("B1").ColumnWidth = 36.42
$ws.Range("H1").ColumnWidth = 61.271
$ws.Range("I1").ColumnWidth = 56.558
